# issue #5: stock data from json to db
#
# The stock ("股票") worksheet gains a `category` column (I) and three new
# trailing columns: `legislator_id` (L), `source_file` (M) and `index` (N).
# The existing date / legislator_name / legislator_id columns (I, J, K)
# shift one column to the right (-> J, K, L).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# ---- Header row (row 1) --------------------------------------------------
# Grab the bold/bordered header style from an existing header cell (K1) and
# stamp it onto the new header cells before writing their text, then shift
# the existing I1/J1/K1 headers right by one column.
$ws.Cells.Item(1, 11).Copy($ws.Cells.Item(1, 12))  # K1 -> L1 (legislator_id)
$ws.Cells.Item(1, 10).Copy($ws.Cells.Item(1, 11))  # J1 -> K1 (legislator_name)
$ws.Cells.Item(1, 9).Copy($ws.Cells.Item(1, 10))   # I1 -> J1 (date)

$ws.Cells.Item(1, 11).Copy($ws.Cells.Item(1, 13))  # style donor -> M1
$ws.Cells.Item(1, 11).Copy($ws.Cells.Item(1, 14))  # style donor -> N1

$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# ---- Data rows (rows 2-8) ------------------------------------------------
for ($r = 2; $r -le 8; $r++) {
    $recordIndex = $ws.Cells.Item($r, 1).Value2   # A: record id (stays the same)

    # Shift existing date / legislator_name / legislator_id right by one
    # column, right-to-left so nothing is clobbered before it's copied.
    $ws.Cells.Item($r, 11).Copy($ws.Cells.Item($r, 12))  # K (legislator_id) -> L
    $ws.Cells.Item($r, 10).Copy($ws.Cells.Item($r, 11))  # J (legislator_name) -> K
    $ws.Cells.Item($r, 9).Copy($ws.Cells.Item($r, 10))   # I (date) -> J

    $ws.Cells.Item($r, 9).Value = "normal"        # I: category
    $ws.Cells.Item($r, 13).Value = "tmp7f9c1"     # M: source_file
    $ws.Cells.Item($r, 14).Value = $recordIndex   # N: index
}
